# Apply the Extent Report update:
# - Insert a new worksheet ("Sheet3" physically, tab named per Excel's
#   auto-numbering) between Sheet1 and Sheet2, containing a new/updated
#   test data table (the "latest Extent Report").
# - Old "Sheet2" keeps its data/name, just shifts to the 3rd tab position.
# - Sheet1's selection moves to column H; the new sheet becomes active.

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet2 = $wb.Worksheets.Item("Sheet2")

# Insert the new worksheet immediately before the existing "Sheet2" tab,
# i.e. right after "Sheet1".
$newSheet = $wb.Worksheets.Add($sheet2)

$cols = @('A','B','C','D','E','F','G','H','I','J','K','L')

$rows = @(
    ,@('Srno', 'Module', 'PageName', 'RunStatus', 'PropertyName', 'PropertyValue', 'Datafield', 'Action', 'Action_Type', 'Test_Case', 'Description', 'Scenario_ID')
    ,@("'1", 'IshinePortal', 'IshineLoginPage', 'Y', $null, $null, $null, 'STARTBROWSER', $null, 'TC_01_01', 'Redirects To Ishine Login page', 'SC_01')
    ,@("'2", 'IshinePortal', 'IshineLoginPage', 'Y', $null, $null, 'URL1', 'BROWSERURL', $null, 'TC_01_01', 'Redirects To Ishine Login page', 'SC_02')
    ,@("'3", 'IshinePortal', 'IshineLoginPage', 'Y', $null, $null, $null, 'wait(2000)', $null, 'TC_01_01', 'Redirects To Ishine Login page', 'SC_03')
    ,@("'4", 'IshinePortal', 'IshineLoginPage', 'Y', 'xpath', '//input[@placeholder="Enter Username"]', 'verify', 'CheckVisibility', 'Element Present or not', 'TC_01_01', 'Redirects To Ishine Login page', 'SC_04')
    ,@("'5", 'IshinePortal', 'IshineLoginPage', 'Y', 'xpath', '//input[@placeholder="Enter Username"]', 'username', 'sendkeys', $null, 'TC_01_01', 'Redirects To Ishine Login page', 'SC_05')
    ,@("'6", 'IshinePortal', 'IshineLoginPage', 'Y', $null, $null, $null, 'wait(1000)', $null, 'TC_01_02', 'User should be able to login after entering credentials.', 'SC_06')
    ,@("'7", 'IshinePortal', 'IshineLoginPage', 'Y', 'xpath', '//input[@placeholder="Enter Password"]', 'password', 'sendkeys', $null, 'TC_01_02', 'User should be able to login after entering credentials.', 'SC_07')
    ,@("'8", 'IshinePortal', 'IshineLoginPage', 'Y', $null, $null, $null, 'wait(1000)', $null, 'TC_01_02', 'User should be able to login after entering credentials.', 'SC_08')
    ,@("'9", 'IshinePortal', 'IshineLoginPage', 'Y', 'xpath', '//button[@type="submit"]', 'verify', 'CheckVisibility', $null, 'TC_01_02', 'User should be able to login after entering credentials.', 'SC_09')
    ,@("'10", 'IshinePortal', 'IshineLoginPage', 'Y', 'xpath', '//button[@type="submit"]', $null, 'click', $null, 'TC_01_02', 'User should be able to login after entering credentials.', 'SC_10')
    ,@("'11", 'IshinePortal', 'MailLoginPage', 'Y', $null, $null, $null, 'wait(4000)', $null, 'TC_01_03', 'Redirects To Apmosys Mail Login page', 'SC_11')
    ,@("'12", 'IshinePortal', 'MailLoginPage', 'Y', $null, $null, $null, 'NewTabOpen', $null, 'TC_01_03', 'Redirects To Apmosys Mail Login page', 'SC_12')
    ,@("'13", 'IshinePortal', 'MailLoginPage', 'Y', $null, $null, $null, 'WindowHandelByIndex(1)', 'Handel the Window', 'TC_01_03', 'Redirects To Apmosys Mail Login page', 'SC_13')
    ,@("'14", 'IshinePortal', 'MailLoginPage', 'Y', $null, $null, 'URL', 'BROWSERURL', $null, 'TC_01_03', 'Redirects To Apmosys Mail Login page', 'SC_14')
    ,@("'15", 'IshinePortal', 'MailLoginPage', 'Y', $null, $null, $null, 'wait(5000)', $null, 'TC_01_03', 'Redirects To Apmosys Mail Login page', 'SC_15')
    ,@("'16", 'IshinePortal', 'MailLoginPage', 'Y', 'xpath', 'email-address', 'verify', 'CheckVisibility', $null, 'TC_01_03', 'Redirects To Apmosys Mail Login page', 'SC_16')
    ,@("'17", 'IshinePortal', 'MailLoginPage', 'Y', 'name', 'email-address', 'username', 'SendKeys', $null, 'TC_01_04', 'User should be able to Go To mailBody And get The OTP', 'SC_17')
    ,@("'18", 'IshinePortal', 'MailLoginPage', 'Y', $null, $null, $null, 'wait(2000)', $null, 'TC_01_04', 'User should be able to Go To mailBody And get The OTP', 'SC_18')
    ,@("'19", 'IshinePortal', 'MailLoginPage', 'Y', 'xpath', '//span[text()="Next"]', $null, 'click', $null, 'TC_01_04', 'User should be able to Go To mailBody And get The OTP', 'SC_19')
    ,@("'20", 'IshinePortal', 'MailLoginPage', 'Y', $null, $null, $null, 'wait(2000)', $null, 'TC_01_04', 'User should be able to Go To mailBody And get The OTP', 'SC_20')
    ,@("'21", 'IshinePortal', 'MailLoginPage', 'Y', 'xpath', '//input[@type="password"]', 'password', 'SendKeys', $null, 'TC_01_04', 'User should be able to Go To mailBody And get The OTP', 'SC_21')
    ,@("'22", 'IshinePortal', 'MailLoginPage', 'Y', $null, $null, $null, 'wait(2000)', $null, 'TC_01_04', 'User should be able to Go To mailBody And get The OTP', 'SC_22')
    ,@("'23", 'IshinePortal', 'MailLoginPage', 'Y', 'name', 'next', $null, 'click', $null, 'TC_01_04', 'User should be able to Go To mailBody And get The OTP', 'SC_23')
    ,@("'24", 'IshinePortal', 'MailLoginPage', 'Y', $null, $null, $null, 'wait(9000)', $null, 'TC_01_04', 'User should be able to Go To mailBody And get The OTP', 'SC_24')
    ,@("'25", 'IshinePortal', 'MailBodyPart', 'Y', 'tagName', 'iframe', $null, 'FRAMECOUNT', $null, 'TC_01_04', 'User should be able to Go To mailBody And get The OTP', 'SC_25')
    ,@("'26", 'IshinePortal', 'MailBodyPart', 'Y', $null, $null, $null, 'wait(3000)', $null, 'TC_01_04', 'User should be able to Go To mailBody And get The OTP', 'SC_26')
    ,@("'27", 'IshinePortal', 'MailBodyPart', 'Y', 'id', 'gui.frm_main.main.mailview#frame', $null, 'FRAMELOCATOR', $null, 'TC_01_04', 'User should be able to Go To mailBody And get The OTP', 'SC_27')
    ,@("'28", 'IshinePortal', 'MailBodyPart', 'Y', $null, $null, $null, 'wait(2000)', $null, 'TC_01_04', 'User should be able to Go To mailBody And get The OTP', 'SC_28')
    ,@("'29", 'IshinePortal', 'MailBodyPart', 'Y', 'xpath', 'iw_webmail_msg_body', 'verify', 'CheckVisibility', $null, 'TC_01_04', 'User should be able to Go To mailBody And get The OTP', 'SC_29')
    ,@("'30", 'IshinePortal', 'MailBodyPart', 'Y', 'id', 'iw_webmail_msg_body', $null, 'gettext', $null, 'TC_01_04', 'User should be able to Go To mailBody And get The OTP', 'SC_30')
    ,@("'31", 'IshinePortal', 'MailBodyPart', 'Y', $null, $null, $null, 'DEFAULTCONTENT', $null, 'TC_01_04', 'User should be able to Go To mailBody And get The OTP', 'SC_31')
    ,@("'32", 'IshinePortal', 'MailBodyPart', 'Y', $null, $null, $null, 'close', $null, 'TC_01_04', 'User should be able to Go To mailBody And get The OTP', 'SC_32')
    ,@("'33", 'IshinePortal', 'MailBodyPart', 'Y', $null, $null, $null, 'wait(3000)', $null, 'TC_01_04', 'User should be able to Go To mailBody And get The OTP', 'SC_33')
    ,@("'34", 'IshinePortal', 'IshineOTPField', 'Y', $null, $null, $null, 'WindowHandelByIndex(0)', $null, 'TC_01_05', 'User should be able to login after entering OTP', 'SC_34')
    ,@("'35", 'IshinePortal', 'IshineOTPField', 'Y', $null, $null, $null, 'wait(10000)', $null, 'TC_01_05', 'User should be able to login after entering OTP', 'SC_35')
    ,@("'36", 'IshinePortal', 'IshineOTPField', 'Y', 'xpath', '//input[@placeholder="Enter OTP"]', $null, 'GetIshineOTP', $null, 'TC_01_05', 'User should be able to login after entering OTP', 'SC_36')
    ,@("'37", 'IshinePortal', 'IshineOTPField', 'Y', 'xpath', '//button[@type="submit"]', 'verify', 'CheckVisibility', $null, 'TC_01_05', 'User should be able to login after entering OTP', 'SC_37')
    ,@("'38", 'IshinePortal', 'IshineOTPField', 'Y', 'xpath', '//button[@type="submit"]', $null, 'click', $null, 'TC_01_05', 'User should be able to login after entering OTP', 'SC_38')
    ,@("'39", 'IshinePortal', 'IshineOTPField', 'Y', $null, $null, $null, 'wait(5000)', $null, 'TC_01_05', 'User should be able to login after entering OTP', 'SC_39')
    ,@('40', 'IshinePortal', 'IshineOTPField', 'Y', $null, $null, $null, 'quit', $null, 'TC_01_05', 'User should be able to login after entering OTP', 'SC_40')
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $rowNum = $i + 1
    $row = $rows[$i]
    for ($j = 0; $j -lt $row.Count; $j++) {
        $val = $row[$j]
        if ($val -ne $null) {
            $cellRef = $cols[$j] + $rowNum
            if ($j -eq 0 -and $rowNum -gt 1 -and -not ($val.StartsWith("'"))) {
                # Row 41's "Srno" is a genuine number, not quote-prefixed text.
                $newSheet.Range($cellRef).Value = [double]$val
            } else {
                $newSheet.Range($cellRef).Value = $val
            }
        }
    }
}

# A handful of "Datafield"/locator cells use the small Consolas font style
# already present in the workbook (from Sheet1). Copy that exact format
# (including the two cells that are styled but left blank).
$styleSourceCell = $sheet1.Range("F17")
$styleSourceCell.Copy()
$consolasCells = @('F17','F18','F24','F28','F29','F30','F31','F36','F37')
foreach ($cellRef in $consolasCells) {
    $newSheet.Range($cellRef).PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

# Sheet1: clear the old "jump to row 22 / F46 selected" view, select column H.
$sheet1.Columns("H").Select()

# New sheet becomes the active tab, selection parked at E30 (matches the
# author's last saved cursor position).
$newSheet.Range("E30").Select()

Write-Host "Inserted sheet:" $newSheet.Name
